# se modifica data para SmokeQa 05-07-2021 R34

$wb = $excel.ActiveWorkbook

# --- Sheet "DatosCuenta" ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokQAJuneLastThree"
$wsCuenta.Range("B2").Value = "SmokeNameQAJuneLastThree"
$wsCuenta.Range("C2").Value = 27100134
$wsCuenta.Range("D2").Value = 135

# --- Sheet "DatosHogar" ---
$wsHogar = $wb.Worksheets.Item("DatosHogar")
$wsHogar.Range("A2").Value = 654

# --- Sheet "DatosMotor" ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Range("A2").Value = "SMP036"
$wsMotor.Range("B2").Value = "ABC12SSMP036"
$wsMotor.Range("C2").Value = "ZAZ123SSMP036"

# --- Sheet "DatosAP" ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Activate()
$wsAP.Range("A2").Value = 21200135
$wsAP.Range("D12").Select()
